$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content fix -----------------------------------------------------
# Row 4 (A4="FilesTab") holds the Cypher query for the Files tab in B4.
# It filtered on experimental_strategies: ["RNA-Seq"]; the input-file
# correction clears that filter back to an empty list so the query
# returns files for every experimental strategy.
$cell = $ws.Range("B4")
$query = $cell.Value()
$fixedQuery = $query.Replace('experimental_strategies: ["RNA-Seq"],', 'experimental_strategies: [],')
$cell.Value = $fixedQuery

# The B/C query columns use a wrapped-text style, so editing the long
# text makes the grid re-autofit row 4's height. Those rows were
# already pinned at Excel's real maximum row height (409.5pt); restore
# that ceiling on rows 2-4 so the edit doesn't blow the row out.
$ws.Rows.Item(2).RowHeight = 409.5
$ws.Rows.Item(3).RowHeight = 409.5
$ws.Rows.Item(4).RowHeight = 409.5

# --- View state --------------------------------------------------------
# Reflect the author scrolling down to / selecting the FilesTab row
# (C4) before saving.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C4").Select()
